# "Add files via upload" - update the projected inventory objective for
# Avril 2025 (C2) from blank to 8000, matching the pattern already used in
# the rows below, and leave the selection on the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 8000
$ws.Range("C2").Select()
